$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Changed default delimiter character from '~' to '.' in header names,
# and dropped the numeric suffix from the vib_wavenumber columns.
$ws.Range("C1").Value = "elements.H"
$ws.Range("D1").Value = "elements.O"
$ws.Range("E1").Value = "elements.Pt"
$ws.Range("L1").Value = "vib_wavenumber"
$ws.Range("M1").Value = "vib_wavenumber"
$ws.Range("N1").Value = "vib_wavenumber"

$ws.Range("C2").Select()
